$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.811.12"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.96%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.911.90"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -2.42%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "587.04"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.54%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.34"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.06%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.505"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.912.12"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -2.38%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -6.42%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +3.54%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -3.82%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000237"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.03%  "
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -2.77%  "
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.72%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.395.12"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -2.31%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.839.13"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.69%  "
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -2.82%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.913.70"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -2.30%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "436.42"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -1.92%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.43"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -1.42%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.660"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -2.84%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -3.41%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "80.89"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -1.82%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.94"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -1.70%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.22"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.02%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -4.50%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0000106"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +19.57%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.20"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.40%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -2.78%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -2.22%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.01%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "25.88"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -3.08%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.976"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.97%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.06"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +2.74%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.50"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -3.44%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -2.62%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -3.55%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.115"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -3.04%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -4.98%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "38.93"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.92%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.697.99"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.67%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "134.97"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.10%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -2.52%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "343.80"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -8.21%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.05%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -1.90%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -4.96%  "
